$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Shape "TekstSylinder 4" - "Clicking a point, updates something"
#   -> two lines: "Clicking a point selects a site " / "(for the second page)"
#   also reflowed/repositioned (spAutoFit textbox, wrap="none")
$shape1 = $s1.Shapes.Item(3)
$tr1 = $shape1.TextFrame.TextRange
$para1 = $tr1.Paragraphs(1, 1)
$chars1 = $para1.Characters(1, $para1.Length)
$chars1.Text = "Clicking a point selects a site "
$tr1.InsertAfter("`r(for the second page)") | Out-Null
$shape1.Left = 43.94015698031496
$shape1.Top = 58.78535273070867
$shape1.Width = 233.04378512755906

# Shape "Rektangel 8" - "Histogram of chosen attribute"
#   -> two lines: "Histogram or plot with time axis" / "of chosen attribute"
$shape2 = $s1.Shapes.Item(6)
$tr2 = $shape2.TextFrame.TextRange
$para2 = $tr2.Paragraphs(1, 1)
$chars2 = $para2.Characters(1, $para2.Length)
$chars2.Text = "Histogram or plot with time axis"
$tr2.InsertAfter("`rof chosen attribute") | Out-Null

# ---------------------------------------------------------------------------
# Slide 2
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Shape "Rektangel 5" - "Weather type" -> "Weather type selector"
$shape3 = $s2.Shapes.Item(4)
$tr3 = $shape3.TextFrame.TextRange
$para3 = $tr3.Paragraphs(1, 1)
$chars3 = $para3.Characters(1, $para3.Length)
$chars3.Text = "Weather type selector"

# Shape "Avrundet rektangel 10" - ARIMA -> ARIMAX / exogeneous -> exogenous
$shape4 = $s2.Shapes.Item(9)
$tr4 = $shape4.TextFrame.TextRange
$para4 = $tr4.Paragraphs(1, 1)
$chars4 = $para4.Characters(1, $para4.Length)
$chars4.Text = "ARIMAX model predicting chosen lice type with weather as exogenous data."
